$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 52
$ws.Range("F3").Value = 965
$ws.Range("F5").Value = 10835
$ws.Range("F6").Value = 1145
$ws.Range("F7").Value = 355
$ws.Range("F8").Value = 589
$ws.Range("F9").Value = 1931
$ws.Range("F10").Value = 586
$ws.Range("F11").Value = 728
$ws.Range("F12").Value = 225
$ws.Range("F13").Value = 299
$ws.Range("F14").Value = 269
$ws.Range("F15").Value = 273
$ws.Range("F16").Value = 987
$ws.Range("F17").Value = 374
$ws.Range("F18").Value = 195
$ws.Range("F19").Value = 414
$ws.Range("F20").Value = 632
$ws.Range("F21").Value = 782
$ws.Range("F22").Value = 175
$ws.Range("F23").Value = 442
$ws.Range("F24").Value = 185

$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 20
$ws.Range("F7").Value = 626

$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 52
$ws.Range("F5").Value = 965
$ws.Range("F8").Value = 10835
$ws.Range("F9").Value = 1145
$ws.Range("F11").Value = 355
$ws.Range("F12").Value = 589
$ws.Range("F13").Value = 1931
$ws.Range("F14").Value = 586
$ws.Range("F15").Value = 728
$ws.Range("F16").Value = 20
$ws.Range("F17").Value = 225
$ws.Range("F18").Value = 299
$ws.Range("F19").Value = 269
$ws.Range("F20").Value = 273
$ws.Range("F21").Value = 987
$ws.Range("F22").Value = 374
$ws.Range("F23").Value = 626
$ws.Range("F24").Value = 195
$ws.Range("F25").Value = 414
$ws.Range("F26").Value = 632
$ws.Range("F27").Value = 782
$ws.Range("F29").Value = 175
$ws.Range("F30").Value = 442
$ws.Range("F31").Value = 185
